$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Payment type" column (E2:E9) to "Done" for every pending entry
$ws.Range("E2:E9").Value = "Done"

# Carry the formatting used across E2:E8 down onto E9 (matches how the
# original edit was made by filling the formula/format down the column)
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the active selection used when making this edit
$ws.Range("E2:E9").Select()
